$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header formatting ------------------------------------------------
# Merge the new Iteration_3/4/5 label cells first (matching the existing
# Iteration_1/Iteration_2 3-column blocks), THEN copy the per-cell style
# of the existing header cells (row 1 / row 2, columns H/I/J) onto the
# new header cells one column at a time. Doing the merge first and the
# style-copy after (rather than merging after formatting) keeps every
# cell on the existing style index (style "1") instead of Excel's
# multi-cell "box border" splitting that merge-after-paste would trigger.
$ws.Range("K1:M1").Merge()
$ws.Range("N1:P1").Merge()
$ws.Range("Q1:S1").Merge()

$ws.Range("H1").Copy()
$ws.Range("K1").PasteSpecial(-4122)
$ws.Range("N1").PasteSpecial(-4122)
$ws.Range("Q1").PasteSpecial(-4122)

$ws.Range("I1").Copy()
$ws.Range("L1").PasteSpecial(-4122)
$ws.Range("O1").PasteSpecial(-4122)
$ws.Range("R1").PasteSpecial(-4122)

$ws.Range("J1").Copy()
$ws.Range("M1").PasteSpecial(-4122)
$ws.Range("P1").PasteSpecial(-4122)
$ws.Range("S1").PasteSpecial(-4122)

$ws.Range("H2").Copy()
$ws.Range("K2").PasteSpecial(-4122)
$ws.Range("N2").PasteSpecial(-4122)
$ws.Range("Q2").PasteSpecial(-4122)

$ws.Range("I2").Copy()
$ws.Range("L2").PasteSpecial(-4122)
$ws.Range("O2").PasteSpecial(-4122)
$ws.Range("R2").PasteSpecial(-4122)

$ws.Range("J2").Copy()
$ws.Range("M2").PasteSpecial(-4122)
$ws.Range("P2").PasteSpecial(-4122)
$ws.Range("S2").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# --- Header row 1: Iteration labels ------------------------------------
$ws.Range("K1").Value = "Iteration_3"
$ws.Range("N1").Value = "Iteration_4"
$ws.Range("Q1").Value = "Iteration_5"

# --- Header row 2: Interval labels -------------------------------------
$ws.Range("K2").Value = "2030"
$ws.Range("L2").Value = "2040"
$ws.Range("M2").Value = "2050"
$ws.Range("N2").Value = "2030"
$ws.Range("O2").Value = "2040"
$ws.Range("P2").Value = "2050"
$ws.Range("Q2").Value = "2030"
$ws.Range("R2").Value = "2040"
$ws.Range("S2").Value = "2050"


# --- Data rows 4-15: numeric values for Iteration_3/4/5 --------------
$ws.Range("K4").Value = [double]"35238095.23809672"
$ws.Range("L4").Value = [double]"0"
$ws.Range("M4").Value = [double]"0"
$ws.Range("N4").Value = [double]"35238095.23810001"
$ws.Range("O4").Value = [double]"0"
$ws.Range("P4").Value = [double]"0"
$ws.Range("Q4").Value = [double]"35238095.23809924"
$ws.Range("R4").Value = [double]"0"
$ws.Range("S4").Value = [double]"0"
$ws.Range("K5").Value = [double]"0"
$ws.Range("L5").Value = [double]"0"
$ws.Range("M5").Value = [double]"0"
$ws.Range("N5").Value = [double]"0"
$ws.Range("O5").Value = [double]"0"
$ws.Range("P5").Value = [double]"0"
$ws.Range("Q5").Value = [double]"0"
$ws.Range("R5").Value = [double]"0"
$ws.Range("S5").Value = [double]"0"
$ws.Range("K6").Value = [double]"1.875326958742353e-06"
$ws.Range("L6").Value = [double]"21050082.96890682"
$ws.Range("M6").Value = [double]"27515499.01900796"
$ws.Range("N6").Value = [double]"-2.29368470520229e-06"
$ws.Range("O6").Value = [double]"21030622.5940332"
$ws.Range("P6").Value = [double]"27515499.01900796"
$ws.Range("Q6").Value = [double]"2.445282953989949e-06"
$ws.Range("R6").Value = [double]"21050125.12373882"
$ws.Range("S6").Value = [double]"27515499.0302716"
$ws.Range("K7").Value = [double]"0"
$ws.Range("L7").Value = [double]"0.001279717154111567"
$ws.Range("M7").Value = [double]"524944.0002144425"
$ws.Range("N7").Value = [double]"1.899009467186991e-06"
$ws.Range("O7").Value = [double]"0"
$ws.Range("P7").Value = [double]"524945.5253625842"
$ws.Range("Q7").Value = [double]"1.775074794745918e-07"
$ws.Range("R7").Value = [double]"-2.878535474074976e-05"
$ws.Range("S7").Value = [double]"524953.17614479"
$ws.Range("K8").Value = [double]"0"
$ws.Range("L8").Value = [double]"0"
$ws.Range("M8").Value = [double]"0"
$ws.Range("N8").Value = [double]"0"
$ws.Range("O8").Value = [double]"0"
$ws.Range("P8").Value = [double]"0"
$ws.Range("Q8").Value = [double]"0"
$ws.Range("R8").Value = [double]"0"
$ws.Range("S8").Value = [double]"0"
$ws.Range("K9").Value = [double]"0"
$ws.Range("L9").Value = [double]"0"
$ws.Range("M9").Value = [double]"0"
$ws.Range("N9").Value = [double]"0"
$ws.Range("O9").Value = [double]"0"
$ws.Range("P9").Value = [double]"0"
$ws.Range("Q9").Value = [double]"0"
$ws.Range("R9").Value = [double]"0"
$ws.Range("S9").Value = [double]"0"
$ws.Range("K10").Value = [double]"0"
$ws.Range("L10").Value = [double]"0"
$ws.Range("M10").Value = [double]"0"
$ws.Range("N10").Value = [double]"0"
$ws.Range("O10").Value = [double]"0"
$ws.Range("P10").Value = [double]"0"
$ws.Range("Q10").Value = [double]"0"
$ws.Range("R10").Value = [double]"0"
$ws.Range("S10").Value = [double]"0"
$ws.Range("K11").Value = [double]"0"
$ws.Range("L11").Value = [double]"0"
$ws.Range("M11").Value = [double]"0"
$ws.Range("N11").Value = [double]"0"
$ws.Range("O11").Value = [double]"0"
$ws.Range("P11").Value = [double]"0"
$ws.Range("Q11").Value = [double]"0"
$ws.Range("R11").Value = [double]"0"
$ws.Range("S11").Value = [double]"0"
$ws.Range("K12").Value = [double]"0"
$ws.Range("L12").Value = [double]"0"
$ws.Range("M12").Value = [double]"0"
$ws.Range("N12").Value = [double]"0"
$ws.Range("O12").Value = [double]"0"
$ws.Range("P12").Value = [double]"0"
$ws.Range("Q12").Value = [double]"0"
$ws.Range("R12").Value = [double]"0"
$ws.Range("S12").Value = [double]"0"
$ws.Range("K13").Value = [double]"0"
$ws.Range("L13").Value = [double]"0"
$ws.Range("M13").Value = [double]"0"
$ws.Range("N13").Value = [double]"0"
$ws.Range("O13").Value = [double]"0"
$ws.Range("P13").Value = [double]"0"
$ws.Range("Q13").Value = [double]"0"
$ws.Range("R13").Value = [double]"0"
$ws.Range("S13").Value = [double]"0"
$ws.Range("K14").Value = [double]"0"
$ws.Range("L14").Value = [double]"0"
$ws.Range("M14").Value = [double]"0"
$ws.Range("N14").Value = [double]"0"
$ws.Range("O14").Value = [double]"0"
$ws.Range("P14").Value = [double]"0"
$ws.Range("Q14").Value = [double]"0"
$ws.Range("R14").Value = [double]"0"
$ws.Range("S14").Value = [double]"0"
$ws.Range("K15").Value = [double]"0"
$ws.Range("L15").Value = [double]"0"
$ws.Range("M15").Value = [double]"0"
$ws.Range("N15").Value = [double]"0"
$ws.Range("O15").Value = [double]"0"
$ws.Range("P15").Value = [double]"0"
$ws.Range("Q15").Value = [double]"0"
$ws.Range("R15").Value = [double]"0"
$ws.Range("S15").Value = [double]"0"
